$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-6 from 45243 to 45244
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
